$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2: test123 -> test
$ws.Range("B2").Value = "test"

# Replace all "XT13" token values with "demo" (D2, D3, D5, D6, D8, D9, D10, D11)
# and fill the previously empty D4 with "demo" as well.
$ws.Range("D2").Value = "demo"
$ws.Range("D3").Value = "demo"
$ws.Range("D4").Value = "demo"
$ws.Range("D5").Value = "demo"
$ws.Range("D6").Value = "demo"
$ws.Range("D8").Value = "demo"
$ws.Range("D9").Value = "demo"
$ws.Range("D10").Value = "demo"
$ws.Range("D11").Value = "demo"

# Column B width change (target stored width 17.7109375; engine quantizes
# ColumnWidth to 1/6-character steps, so 16.8 is the closest achievable input)
$ws.Columns.Item(2).ColumnWidth = 16.8

# Update the selected cell/active selection on the sheet
$ws.Range("D11").Select()
